$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Metadata")

# --- Data sheet: insert 3 new rows at the top (2024, 2023, 2022) ---
$ws1.Rows("2:4").Insert()

$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "2024"
$ws1.Range("A2").Style = "Normal"
$ws1.Range("B2").Value = 3.1

$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "2023"
$ws1.Range("A3").Style = "Normal"
$ws1.Range("B3").Value = 3.6

$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = "2022"
$ws1.Range("A4").Style = "Normal"
$ws1.Range("B4").Value = 3.2

# --- Metadata sheet ---
# A1 held a blank/empty shared string; retarget it to the single-space string
# so the now-unused empty string entry drops out of sharedStrings.
$ws2.Range("A1").Value = " "

# Update "observaciones" value (row 8, column B)
$ws2.Range("B8").Value = "Desde marzo de 2020 hasta junio de 2021 se interrumpió el relevamiento presencial y se aplicó de manera telefónica un cuestionario restringido con el objetivo de continuar publicando los indicadores de ingresos y mercado de trabajo. En ese período la encuesta pasó a ser de paneles rotativos elegidos al azar a partir de los casos respondentes del año anterior. En julio de 2021 el INE retomó la realización de encuestas presenciales, pero introdujo un cambio metodológico, ya que la ECH pasa a ser una encuesta de panel rotativo con periodicidad mensual compuesta por seis paneles o grupos de rotación, cada uno de los cuales es una muestra representativa de la población. Con esta nueva metodología, cada hogar seleccionado participa durante seis meses de la ECH. Durante el año 2020 y hasta julio del año 2021 se suspende el relevamiento de la información necesaria para construir indicadores relativos al nivel y la trayectoria educativa. A partir de esta fecha, las preguntas se relevan en el formulario presencial. Un conjunto importante de indicadores educativos tienen un efecto estacional, por lo que no se recomienda comparar los resultados del segundo semestre del 2021 con la información anual. Las estimaciones desde 2022 se calculan a partir de la muestra de implantación. Respecto a la forma de preguntar asistencia a centros educativos, el INE realizó un cambio metodológico en el relevamiento. Anteriormente, se consultaba a las personas por la asistencia a cada nivel educativo. Se generaban, así, ocho variables de asistencia, una correspondiente a cada nivel.  El porcentaje de personas que no asisten era un indicador resumen de esta información.  A partir de 2020, se consulta a las personas si asisten a un establecimiento de enseñanaza de manera general, mediante una única pregunta. A partir de esta fecha, el porcentaje de no asistentes se calcula únicamente a partir de esta pregunta. A su vez, se modifica la forma de relevamiento en la culminación de ciclos educativos. En particular, cambia el relevamiento de cantidad de años aprobados en UTU. Hasta el año 2019 se relevaban los años aprobados en bachillerato tecnológico y en educación técnica. En el segundo caso era posible distinguir el curso según la exigencia previa para cursarlo. A partir de julio de 2021 se distinguen los años de Educación Media Básica y Educación Media Superior, tanto de liceo como de CEPT-UTU. Además, se consulta de forma independiente los años realizados en cursos técnicos en CEPT-UTU. No se releva el nivel de exigencia previa para asistir a estos cursos."

# Insert a new row for "actualizacion" / "Julio 2025" before "cita" (was row 9, now row 9 after insert)
$ws2.Rows("9:9").Insert()
$ws2.Range("A9").Value = "actualizacion"
$ws2.Range("B9").Value = "Julio 2025"

# Update "cita" value (now row 10) with the new citation text (includes trailing newline)
$ws2.Range("B10").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE (Hasta 2019) / A partir de 2020 con base en ECH - INE`n"

Write-Output "done"
